$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hanz Jansen123"
$ws.Range("C2").Value = 123412313

$ws.Range("A3").Value = "batmat"
$ws.Range("B3").Value = "IT DEPARTMENT"
$ws.Range("C3").Value = 1234

Write-Host ("A2: " + $ws.Range("A2").Text)
Write-Host ("C2: " + $ws.Range("C2").Text)
Write-Host ("A3: " + $ws.Range("A3").Text)
Write-Host ("B3: " + $ws.Range("B3").Text)
Write-Host ("C3: " + $ws.Range("C3").Text)
